# Add three new vocabulary rows (hustle, hostile, construe) to the word list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 89: hustle
$ws.Range("A89").Value2 = "hustle"
$ws.Range("B89").Value2 = "to make someone move quickly, especially by pushing them roughly"
$ws.Range("C89").Value2 = "I was hustled out of the building by a couple of security men."
$ws.Range("D89").Value2 = "they hustled to finish the task on time."

# Row 90: hostile (filled word, example2, meaning, example1 - matches source order)
$ws.Range("A90").Value2 = "hostile"
$ws.Range("D90").Value2 = "the boy feels hostile towards his father."
$ws.Range("B90").Value2 = "angry and deliberately unfriendly towards someone, and ready to argue with them"
$ws.Range("C90").Value2 = "Southampton fans gave their former coach a hostile reception."

# Row 91: construe
$ws.Range("A91").Value2 = "construe"
$ws.Range("B91").Value2 = "to comprehend or explain the meaning or intention of; assign a meaning to; interpret."
$ws.Range("C91").Value2 = "comments that could be construed as sexist"
$ws.Range("D91").Value2 = "they construcedmy words as cirtical and hostile."

# Match the row heights used for the newly-added rows.
$ws.Rows(89).RowHeight = 60
$ws.Rows(90).RowHeight = 75
$ws.Rows(91).RowHeight = 75

# Update the view so the new rows are visible/selected, mirroring the authored edit.
$excel.ActiveWindow.ScrollRow = 83
$ws.Range("A89:D91").Select()
